# Weekly update for "Hortaliza, Femacal de La Calera - Achicoria"
# A new weekly record is inserted at row 116 (shifting the existing
# rows 116-143 down to 117-144, so the previous last row, 143, becomes
# the new row 144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 116, pushing rows 116:143 down to 117:144.
# -4121 = xlShiftDown
$ws.Range("A116:R116").Insert(-4121)

# Populate the newly inserted row 116 with the new weekly observation.
$ws.Range("A116").Value = 3
$ws.Range("B116").Value = "Femacal de La Calera"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44511
$ws.Range("E116").Value = 5
$ws.Range("F116").Value = 100112010
$ws.Range("G116").Value = "Achicoria"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 115
$ws.Range("K116").Value = 5800
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 5896
$ws.Range("N116").Value = "`$/caja 16 unidades"
$ws.Range("O116").Value = "Provincia de Quillota"
$ws.Range("P116").Value = 368
$ws.Range("Q116").Value = 16
$ws.Range("R116").Value = "Hortaliza"
